$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.845.13"
$ws.Range("E2").Value = "  +1.44%  "

$ws.Range("D3").Value = "2.663.04"
$ws.Range("E3").Value = "  +3.06%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.62"
$ws.Range("E5").Value = "  +2.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.99"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  -0.89%  "

$ws.Range("E9").Value = "  +0.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.66"
$ws.Range("E10").Value = "  +0.01%  "

$ws.Range("E11").Value = "  -0.02%  "

$ws.Range("E12").Value = "  +0.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.67"
$ws.Range("E13").Value = "  +1.12%  "

$ws.Range("D14").Value = "3.137.52"
$ws.Range("E14").Value = "  +2.93%  "

$ws.Range("D15").Value = "63.771.42"
$ws.Range("E15").Value = "  +1.50%  "

$ws.Range("E16").Value = "  +0.53%  "

$ws.Range("D17").Value = "2.677.35"
$ws.Range("E17").Value = "  +3.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.44"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "343.69"
$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.39"
$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.80"
$ws.Range("E21").Value = "  +1.70%  "

$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.16"
$ws.Range("E23").Value = "  +1.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.64"
$ws.Range("E24").Value = "  +14.03%  "

$ws.Range("E25").Value = "  +5.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "575.86"
$ws.Range("E26").Value = "  +23.56%  "

$ws.Range("E27").Value = "  -0.98%  "

$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.08"
$ws.Range("E28").Value = "  +2.58%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.57"
$ws.Range("E29").Value = "  +2.76%  "

$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.01"
$ws.Range("E31").Value = "  +3.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.80"
$ws.Range("E32").Value = "  +11.94%  "

$ws.Range("E33").Value = "  -0.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "175.31"
$ws.Range("E34").Value = "  +0.17%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.403"
$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.78"
$ws.Range("E37").Value = "  +5.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.23"
$ws.Range("E38").Value = "  +0.61%  "

$ws.Range("E39").Value = "  +3.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "171.26"
$ws.Range("E40").Value = "  +8.05%  "

$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.50"
$ws.Range("E42").Value = "  +2.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.78"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.88"
$ws.Range("E44").Value = "  +2.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.633"
$ws.Range("E45").Value = "  -0.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0554"
$ws.Range("E46").Value = "  +1.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0241"
$ws.Range("E47").Value = "  +1.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0963"
$ws.Range("E48").Value = "  -0.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.80"
$ws.Range("E49").Value = "  +1.73%  "

$ws.Range("E50").Value = "  +2.30%  "

$ws.Range("E51").Value = "  +11.97%  "
